$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: classical-best-embeddings -> classical-best-embed (label), and updated values
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.075
$ws.Range("D2").Value = 0.058
$ws.Range("E2").Value = 0.054
$ws.Range("I2").Value = 0.055
$ws.Range("J2").Value = 0.064

# Row 3: BERT-base vs. classical-best-tfidf (label unchanged), values updated
$ws.Range("C3").Value = 0.068
$ws.Range("D3").Value = 0.154
$ws.Range("E3").Value = 0.153
$ws.Range("F3").Value = 0.122
$ws.Range("G3").Value = 0.167
$ws.Range("H3").Value = 0.18
$ws.Range("I3").Value = 0.124
$ws.Range("J3").Value = 0.141

# Row 4: BERT-base vs. classical-best-embeddings -> BERT-base vs. classical-best-embed (label), values updated
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("C4").Value = -0.007
$ws.Range("D4").Value = 0.096
$ws.Range("E4").Value = 0.099
$ws.Range("F4").Value = 0.08699999999999999
$ws.Range("G4").Value = 0.08799999999999999
$ws.Range("H4").Value = 0.095
$ws.Range("I4").Value = 0.06900000000000001
$ws.Range("J4").Value = 0.076

# Row 5: BERT-base-nli vs. classical-best-tfidf (label unchanged), values updated
$ws.Range("B5").Value = 0.281
$ws.Range("C5").Value = 0.222
$ws.Range("D5").Value = 0.189
$ws.Range("E5").Value = 0.186
$ws.Range("G5").Value = 0.139
$ws.Range("H5").Value = 0.154
$ws.Range("I5").Value = 0.184
$ws.Range("J5").Value = 0.171

# Row 6: BERT-base-nli vs. classical-best-embeddings -> BERT-base-nli vs. classical-best-embed (label), values updated
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.281
$ws.Range("C6").Value = 0.147
$ws.Range("D6").Value = 0.131
$ws.Range("E6").Value = 0.132
$ws.Range("G6").Value = 0.06
$ws.Range("H6").Value = 0.06900000000000001
$ws.Range("I6").Value = 0.128
$ws.Range("J6").Value = 0.107

# Row 7: BERT-base-nli vs. BERT-base (label unchanged), values updated
$ws.Range("B7").Value = 0.281
$ws.Range("C7").Value = 0.154
$ws.Range("D7").Value = 0.035
$ws.Range("E7").Value = 0.033
$ws.Range("F7").Value = 0.016
$ws.Range("G7").Value = -0.028
$ws.Range("H7").Value = -0.026
$ws.Range("I7").Value = 0.06
$ws.Range("J7").Value = 0.031
